$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 40 ---
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Range("A40").Value = "Beschadigd product ontvangen"
$wsLogs.Range("B40").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C40").Value = "Het product dat ik heb ontvangen is beschadigd aangekomen."
$wsLogs.Range("D40").Value = "Overig"
$wsLogs.Range("F40").Value = "2025-06-24 22:00:20"
$wsLogs.Range("G40").Value = "Nee"

# Expand the conditional formatting ranges that used to stop at row 39
$fcsD = $wsLogs.Range("D2:D39").FormatConditions
for ($i = 1; $i -le $fcsD.Count; $i++) {
    $fcsD.Item($i).ModifyAppliesToRange($wsLogs.Range("D2:D40"))
}

$fcsG = $wsLogs.Range("G2:G39").FormatConditions
for ($i = 1; $i -le $fcsG.Count; $i++) {
    $fcsG.Item($i).ModifyAppliesToRange($wsLogs.Range("G2:G40"))
}

# --- Dashboard sheet: append new row 9 (Overig / 1) ---
$wsDash = $wb.Worksheets.Item("Dashboard")

$wsDash.Range("A9").Value = "Overig"
$wsDash.Range("B9").Value = 1

# --- Chart: extend the category / value series ranges from row 8 to row 9 ---
$chartObj = $wsDash.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$9"
$series.Values = "='Dashboard'!`$B`$2:`$B`$9"
